$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI metrics (ligand/receptor expression & specificity
# columns G, H, M-T) for the Fgf5-Fgfr1 LR-pair sheet, rows 2-7.

# Row 2 (Target cluster: ECs)
$ws.Cells.Item(2, 7).Value = 0.4859026666666666
$ws.Cells.Item(2, 8).Value = 1.457708
$ws.Cells.Item(2, 13).Value = 10.48767733333333
$ws.Cells.Item(2, 14).Value = 31.463032
$ws.Cells.Item(2, 15).Value = 0.1222087640673552
$ws.Cells.Item(2, 16).Value = 0.1222087640673552
$ws.Cells.Item(2, 17).Value = 5.095990383406222
$ws.Cells.Item(2, 18).Value = 45.863913450656
$ws.Cells.Item(2, 19).Value = 0.1222087640673552
$ws.Cells.Item(2, 20).Value = 0.1222087640673552

# Row 3 (Target cluster: FAPs)
$ws.Cells.Item(3, 7).Value = 0.4859026666666666
$ws.Cells.Item(3, 8).Value = 1.457708
$ws.Cells.Item(3, 15).Value = 0.7340790765058636
$ws.Cells.Item(3, 16).Value = 0.7340790765058635
$ws.Cells.Item(3, 17).Value = 30.61040624281111
$ws.Cells.Item(3, 18).Value = 275.4936561853
$ws.Cells.Item(3, 19).Value = 0.7340790765058636
$ws.Cells.Item(3, 20).Value = 0.7340790765058635

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Cells.Item(4, 7).Value = 0.4859026666666666
$ws.Cells.Item(4, 8).Value = 1.457708
$ws.Cells.Item(4, 13).Value = 0.3322793333333333
$ws.Cells.Item(4, 14).Value = 0.996838
$ws.Cells.Item(4, 15).Value = 0.003871919907635547
$ws.Cells.Item(4, 16).Value = 0.003871919907635547
$ws.Cells.Item(4, 17).Value = 0.1614554141448889
$ws.Cells.Item(4, 18).Value = 1.453098727304
$ws.Cells.Item(4, 19).Value = 0.003871919907635547
$ws.Cells.Item(4, 20).Value = 0.003871919907635547

# Row 5 (Target cluster: MuSCs)
$ws.Cells.Item(5, 7).Value = 0.4859026666666666
$ws.Cells.Item(5, 8).Value = 1.457708
$ws.Cells.Item(5, 13).Value = 10.25458433333333
$ws.Cells.Item(5, 14).Value = 30.763753
$ws.Cells.Item(5, 15).Value = 0.1194926233493133
$ws.Cells.Item(5, 16).Value = 0.1194926233493133
$ws.Cells.Item(5, 17).Value = 4.982729873124889
$ws.Cells.Item(5, 18).Value = 44.844568858124
$ws.Cells.Item(5, 19).Value = 0.1194926233493133
$ws.Cells.Item(5, 20).Value = 0.1194926233493133

# Row 6 (Target cluster: Neutrophils)
$ws.Cells.Item(6, 7).Value = 0.4859026666666666
$ws.Cells.Item(6, 8).Value = 1.457708
$ws.Cells.Item(6, 13).Value = 0.7572163333333334
$ws.Cells.Item(6, 14).Value = 2.271649
$ws.Cells.Item(6, 15).Value = 0.008823543029319092
$ws.Cells.Item(6, 16).Value = 0.00882354302931909
$ws.Cells.Item(6, 17).Value = 0.3679334356102222
$ws.Cells.Item(6, 18).Value = 3.311400920492
$ws.Cells.Item(6, 19).Value = 0.008823543029319092
$ws.Cells.Item(6, 20).Value = 0.00882354302931909

# Row 7 (Target cluster: Resolving-Mac)
$ws.Cells.Item(7, 7).Value = 0.4859026666666666
$ws.Cells.Item(7, 8).Value = 1.457708
$ws.Cells.Item(7, 13).Value = 0.9889696666666666
$ws.Cells.Item(7, 14).Value = 2.966909
$ws.Cells.Item(7, 15).Value = 0.01152407314051338
$ws.Cells.Item(7, 16).Value = 0.01152407314051338
$ws.Cells.Item(7, 17).Value = 0.4805429982857777
$ws.Cells.Item(7, 18).Value = 4.324886984571999
$ws.Cells.Item(7, 19).Value = 0.01152407314051338
$ws.Cells.Item(7, 20).Value = 0.01152407314051338
